# Regenerate orders with updated distance/sizes.
# The shared-strings table encodes trial condition / filename tokens such as
# "D64", "D80", "D51" (Distance) and "S30" (Size). This run regenerates the
# order with new distance/size values:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These tokens appear (sometimes as a whole cell, sometimes as part of a
# longer filename string) in the Condition, Filename_Left, Filename_Right,
# Distance and Size columns. A sheet-wide substring Find/Replace reproduces
# every occurrence without needing to know each cell address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

# Order matters only in that none of the replacement tokens (D69/D86/D55/S31)
# collide with any of the other search tokens, so a straightforward
# sequential replace is safe.
$rng.Replace("D64", "D69") | Out-Null
$rng.Replace("D80", "D86") | Out-Null
$rng.Replace("D51", "D55") | Out-Null
$rng.Replace("S30", "S31") | Out-Null
